$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object "object[,]" 24,5
$il = New-Object "object[,]" 24,6

$bf[0,0] = 1.02
$bf[0,1] = 1.037535968089784
$bf[0,2] = 1.047018807962933
$bf[0,3] = 1.036314708015164
$bf[0,4] = 1.057583501978786
$il[0,0] = 1.0416363482431
$il[0,1] = 1.042638307426631
$il[0,2] = 1.049782618260352
$il[0,3] = 1.039108802664993
$il[0,4] = 1.060318073686908
$il[0,5] = 1.044118973415462

$bf[1,0] = 1.02
$bf[1,1] = 1.038392959414711
$bf[1,2] = 1.047710356703334
$bf[1,3] = 1.037041227658753
$bf[1,4] = 1.058438122700043
$il[1,0] = 1.041851486440149
$il[1,1] = 1.043140183149734
$il[1,2] = 1.050286238753185
$il[1,3] = 1.039645118768165
$il[1,4] = 1.060986455766684
$il[1,5] = 1.044621561859658

$bf[2,0] = 1.02
$bf[2,1] = 1.038948153058497
$bf[2,2] = 1.048158426837216
$bf[2,3] = 1.037512275436404
$bf[2,4] = 1.058992096763948
$il[2,0] = 1.041989876830475
$il[2,1] = 1.043464932975555
$il[2,2] = 1.050612028034786
$il[2,3] = 1.03999242175089
$il[2,4] = 1.061419274351417
$il[2,5] = 1.044946772867481

$bf[3,0] = 1.02
$bf[3,1] = 1.039181713611396
$bf[3,2] = 1.048346935466272
$bf[3,3] = 1.037710527429864
$bf[3,4] = 1.059225219534524
$il[3,0] = 1.04204785950206
$il[3,1] = 1.043601457198713
$il[3,2] = 1.050748967741826
$il[3,3] = 1.040138491363459
$il[3,4] = 1.061601308850251
$il[3,5] = 1.045083490970695

$bf[4,0] = 1.02
$bf[4,1] = 1.039220938625538
$bf[4,2] = 1.048378595050199
$bf[4,3] = 1.037743827867139
$bf[4,4] = 1.059264375427331
$il[4,0] = 1.042057583495311
$il[4,1] = 1.043624380144983
$il[4,2] = 1.050771959186654
$il[4,3] = 1.040163020775148
$il[4,4] = 1.061631877745569
$il[4,5] = 1.045106446470178

$bf[5,0] = 1.02
$bf[5,1] = 1.038951273290518
$bf[5,2] = 1.048160945151051
$bf[5,3] = 1.037514923613231
$bf[5,4] = 1.058995210852302
$il[5,0] = 1.04199065237136
$il[5,1] = 1.043466757222836
$il[5,2] = 1.050613857918365
$il[5,3] = 1.039994373292618
$il[5,4] = 1.061421706402525
$il[5,5] = 1.044948599705402

$bf[6,0] = 1.02
$bf[6,1] = 1.037825454108851
$bf[6,2] = 1.047252396499702
$bf[6,3] = 1.036560042830785
$bf[6,4] = 1.057872121554011
$il[6,0] = 1.041709224227724
$il[6,1] = 1.04280791737545
$il[6,2] = 1.049952836357088
$il[6,3] = 1.03928999636519
$il[6,4] = 1.060543887106532
$il[6,5] = 1.044288824229862

$bf[7,0] = 1.02
$bf[7,1] = 1.035846763681862
$bf[7,2] = 1.045656027851019
$bf[7,3] = 1.034884705838522
$bf[7,4] = 1.055900663645517
$il[7,0] = 1.041207075746784
$il[7,1] = 1.041647032887566
$il[7,2] = 1.048787424858885
$il[7,3] = 1.038050926146663
$il[7,4] = 1.058999660855606
$il[7,5] = 1.043126291152852

$bf[8,0] = 1.02
$bf[8,1] = 1.034531190472219
$bf[8,2] = 1.044594985931969
$bf[8,3] = 1.033772817929168
$bf[8,4] = 1.054591560452129
$il[8,0] = 1.040868160283168
$il[8,1] = 1.040873234814864
$il[8,2] = 1.04801015411057
$il[8,3] = 1.037226392881189
$il[8,4] = 1.05797202267742
$il[8,5] = 1.04235139419813

$bf[9,0] = 1.02
$bf[9,1] = 1.033962395242903
$bf[9,2] = 1.044136324996521
$bf[9,3] = 1.033292567218569
$bf[9,4] = 1.054025962066183
$il[9,0] = 1.040720430923833
$il[9,1] = 1.040538216961669
$il[9,2] = 1.047673525614519
$il[9,3] = 1.036869736782578
$il[9,4] = 1.057527501999932
$il[9,5] = 1.042015900581156

$bf[10,0] = 1.02
$bf[10,1] = 1.033751249638818
$bf[10,2] = 1.0439660764531
$bf[10,3] = 1.033114363431529
$bf[10,4] = 1.053816063663605
$il[10,0] = 1.040665411635363
$il[10,1] = 1.040413783950755
$il[10,2] = 1.047548478344142
$il[10,3] = 1.036737316117137
$il[10,4] = 1.057362456974898
$il[10,5] = 1.041891290861094

$bf[11,0] = 1.02
$bf[11,1] = 1.033796535167457
$bf[11,2] = 1.044002589950356
$bf[11,3] = 1.033152580451734
$bf[11,4] = 1.053861078953586
$il[11,0] = 1.040677220056016
$il[11,1] = 1.040440474892816
$il[11,2] = 1.047575301769576
$il[11,3] = 1.03676571818501
$il[11,4] = 1.057397856510759
$il[11,5] = 1.041918019707355

$bf[12,0] = 1.02
$bf[12,1] = 1.033944939214858
$bf[12,2] = 1.044122249770129
$bf[12,3] = 1.033277833100256
$bf[12,4] = 1.054008607911195
$il[12,0] = 1.040715885987039
$il[12,1] = 1.040527931129174
$il[12,2] = 1.047663189333015
$il[12,3] = 1.036858789672483
$il[12,4] = 1.057513857894469
$il[12,5] = 1.042005600141599

$bf[13,0] = 1.02
$bf[13,1] = 1.034036393155374
$bf[13,2] = 1.044195991896654
$bf[13,3] = 1.033355029655807
$bf[13,4] = 1.054099530608514
$il[13,0] = 1.040739690014759
$il[13,1] = 1.040581816838982
$il[13,2] = 1.047717338680372
$il[13,3] = 1.0369161417199
$il[13,4] = 1.057585339466217
$il[13,5] = 1.042059562375296

$bf[14,0] = 1.02
$bf[14,1] = 1.034568957720673
$bf[14,2] = 1.044625442274878
$bf[14,3] = 1.033804716110425
$bf[14,4] = 1.054629123895994
$il[14,0] = 1.040877944076711
$il[14,1] = 1.040895469834486
$il[14,2] = 1.048032493776164
$il[14,3] = 1.037250070931061
$il[14,4] = 1.058001533742692
$il[14,5] = 1.04237366079403

$bf[15,0] = 1.02
$bf[15,1] = 1.034903251857496
$bf[15,2] = 1.044895034492419
$bf[15,3] = 1.034087116189707
$bf[15,4] = 1.054961660315877
$il[15,0] = 1.040964406203886
$il[15,1] = 1.041092228152676
$il[15,2] = 1.04823016573998
$il[15,3] = 1.037459636526676
$il[15,4] = 1.058262723959696
$il[15,5] = 1.0425706985316

$bf[16,0] = 1.02
$bf[16,1] = 1.035098322520837
$bf[16,2] = 1.045052357942527
$bf[16,3] = 1.034251951434253
$bf[16,4] = 1.055155743847531
$il[16,0] = 1.041014743811735
$il[16,1] = 1.041206997901547
$il[16,2] = 1.048345458078719
$il[16,3] = 1.037581908417947
$il[16,4] = 1.058415115435261
$il[16,5] = 1.042685631266679

$bf[17,0] = 1.02
$bf[17,1] = 1.0351648505301
$bf[17,2] = 1.045106013818571
$bf[17,3] = 1.034308175622999
$bf[17,4] = 1.055221941724199
$il[17,0] = 1.041031891620541
$il[17,1] = 1.041246132051987
$il[17,2] = 1.048384768665411
$il[17,3] = 1.037623606008115
$il[17,4] = 1.058467084343175
$il[17,5] = 1.042724820992101

$bf[18,0] = 1.02
$bf[18,1] = 1.034867376710586
$bf[18,2] = 1.044866102041784
$bf[18,3] = 1.034056805316528
$bf[18,4] = 1.054925969810334
$il[18,0] = 1.040955139383413
$il[18,1] = 1.041071117429545
$il[18,2] = 1.04820895805867
$il[18,3] = 1.037437148406196
$il[18,4] = 1.058234696195588
$il[18,5] = 1.042549557828821

$bf[19,0] = 1.02
$bf[19,1] = 1.033901234313081
$bf[19,2] = 1.044087009648786
$bf[19,3] = 1.033240944258894
$bf[19,4] = 1.053965159049205
$il[19,0] = 1.040704503859372
$il[19,1] = 1.040502177229241
$il[19,2] = 1.047637308858959
$il[19,3] = 1.036831380846739
$il[19,4] = 1.057479696439827
$il[19,5] = 1.041979809668174

$bf[20,0] = 1.02
$bf[20,1] = 1.033294536721066
$bf[20,2] = 1.043597850113668
$bf[20,3] = 1.032729038028927
$bf[20,4] = 1.053362159550518
$il[20,0] = 1.040546074847752
$il[20,1] = 1.040144506168927
$il[20,2] = 1.047277842037308
$il[20,3] = 1.036450842874913
$il[20,4] = 1.057005402465812
$il[20,5] = 1.041621630673929

$bf[21,0] = 1.02
$bf[21,1] = 1.033616086553217
$bf[21,2] = 1.043857097145415
$bf[21,3] = 1.033000308317322
$bf[21,4] = 1.05368171604697
$il[21,0] = 1.04063014087466
$il[21,1] = 1.040334109754479
$il[21,2] = 1.047468406333972
$il[21,3] = 1.036652541315753
$il[21,4] = 1.057256795806408
$il[21,5] = 1.041811503518321

$bf[22,0] = 1.02
$bf[22,1] = 1.034883586890368
$bf[22,2] = 1.044879175139733
$bf[22,3] = 1.034070501137041
$bf[22,4] = 1.054942096440921
$il[22,0] = 1.04095932695212
$il[22,1] = 1.041080656444574
$il[22,2] = 1.048218540916837
$il[22,3] = 1.037447309709329
$il[22,4] = 1.058247360601106
$il[22,5] = 1.042559110390346

$bf[23,0] = 1.02
$bf[23,1] = 1.036357683062707
$bf[23,2] = 1.046068170341111
$bf[23,3] = 1.035316946437178
$bf[23,4] = 1.056409424157513
$il[23,0] = 1.041337627755927
$il[23,1] = 1.041947132769334
$il[23,2] = 1.049088775530134
$il[23,3] = 1.038370994595065
$il[23,4] = 1.059398561987741
$il[23,5] = 1.043426817210873

$ws.Range("B2:F25").Value = $bf
$ws.Range("I2:N25").Value = $il

